$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: fill in the Time Clocked-Out and Total Duration for the 2026-01-25 entry
$ws.Range("C4").Value = "18:28:07"
$ws.Range("D4").Value = "4.86 Hours"

# Row 5: a second clock-in on 2026-01-25, plus the overall Total Duration label/value
# Leading apostrophe keeps the date/time strings as literal text instead of
# being auto-converted to date/time serial numbers.
$ws.Range("A5").Value = "'2026-01-25"
$ws.Range("B5").Value = "18:58:12"
$ws.Range("C5").Value = "Total Duration:"
$ws.Range("D5").Value = "5 Hours"

# The cells above were brand new (previously blank), so they picked up the
# default style. Copy the formatting already used throughout the table
# (e.g. C6, which shares the row's style) onto the new/edited cells so they
# stay visually consistent with the rest of the sheet.
$ws.Range("C6").Copy()
$ws.Range("C4:D5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A5:B5").PasteSpecial(-4122)  # xlPasteFormats
